$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J header and per-row "Mitigacao de ilhas de calor" values, mirroring the
# existing Sim/Nao columns (F-I) that were added alongside it.
$jValues = @("Não","Não","Não","Não","Não","Sim","Não","Não","Não","Não","Não","Sim","Sim","Sim","Sim","Sim","Não","Sim","Não","Sim","Não","Sim","Sim","Sim","Sim","Sim","Sim","Sim","Sim","Sim","Sim","Sim","Não","Sim","Sim","Não","Não","Sim","Sim","Sim","Não","Não","Sim","Sim","Sim","Não","Sim","Sim","Sim","Não","Sim","Sim","Não","Sim","Sim","Não","Sim","Sim","Sim","Sim","Sim","Sim","Não","Não","Não","Não","Não","Sim","Não","Não","Sim")

# Copy formatting from column I (last existing data column) into the new column J so the
# new cells inherit the same cell style used by the rest of the table.
$ws.Range("I1:I72").Copy()
$ws.Range("J1:J72").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 10).Value = "Mitigação de ilhas de calor"

for ($i = 0; $i -lt $jValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}

# Give the new column a sensible width (closest this engine's column-width quantization
# allows to the authored 28.44140625 character-width value).
$ws.Columns.Item(10).ColumnWidth = 27.65

# Re-apply the AutoFilter so it spans the new column, then refresh the hidden
# _FilterDatabase defined name that Excel keeps in sync with the filter range.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:J72").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "areas_verdes_cau_dgat_sqa!_FilterDatabase") {
        $n.RefersTo = "=areas_verdes_cau_dgat_sqa!`$A`$1:`$J`$72"
    }
}

# Move the view so the new column is visible and selected, matching the saved view state.
$ws.Range("E1").Select()
$ws.Range("J1").Select()
